# Update countries & provincias Spain
#
# The source "Pais" sheet lists one country per row (rows 4-219), sorted
# descending by total cases. This refresh (26-Jun-2020, 14:34 -> 15:51)
# brings in newer case counts for a handful of countries; a few of those
# updates cause countries that were tied (or become tied) on total cases
# to swap display order, which is why some rows only get a new country
# name in column A while their case numbers stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell: refreshed timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 15:51"

# --- Country name swaps caused by rows trading rank order ---
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("A80").Value = "Tayikistan"
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Laos"
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Papua Nueva Guinea"
$ws.Range("A213").Value = "Montserrat"

# --- Updated case figures (Casos totales / Nuevos casos / Casos activos /
#     Recuperados / Casos criticos / Muertes hoy / Muertes) ---

# Estados Unidos
$ws.Range("B4").Value = 2505909
$ws.Range("C4").Value = 1321
$ws.Range("D4").Value = 1052442
$ws.Range("E4").Value = 1326644
$ws.Range("G4").Value = 43
$ws.Range("H4").Value = 126823

# Arabia Saudita
$ws.Range("B18").Value = 174577
$ws.Range("C18").Value = 3938
$ws.Range("D18").Value = 120471
$ws.Range("E18").Value = 52632
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = 1474

# row 34
$ws.Range("B34").Value = 46973
$ws.Range("C34").Value = 410
$ws.Range("D34").Value = 35469
$ws.Range("E34").Value = 11194
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 310

# Emiratos Arabes Unidos
$ws.Range("B38").Value = 40866
$ws.Range("C38").Value = 451
$ws.Range("D38").Value = 26633
$ws.Range("E38").Value = 12678
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = 1555

# row 62
$ws.Range("B62").Value = 13565
$ws.Range("C62").Value = 193
$ws.Range("D62").Value = 12232
$ws.Range("E62").Value = 1068
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 265

# Finlandia
$ws.Range("E76").Value = 263
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 328

# Now "Republica de Macedonia" (A79)
$ws.Range("B79").Value = 5758
$ws.Range("C79").Value = 163
$ws.Range("D79").Value = 2206
$ws.Range("E79").Value = 3284
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 268

# Now "Tayikistan" (A80)
$ws.Range("B80").Value = 5691
$ws.Range("D80").Value = 4267
$ws.Range("E80").Value = 1372
$ws.Range("H80").Value = 52

# Islandia
$ws.Range("B113").Value = 1832
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 1813

# Libano
$ws.Range("B115").Value = 1697
$ws.Range("C115").Value = 35
$ws.Range("E115").Value = 520

# Now "Papua Nueva Guinea" (A211)
$ws.Range("C211").Value = 1
$ws.Range("D211").Value = 8
$ws.Range("E211").Value = 3
$ws.Range("H211").Value = 0

# Now "Montserrat" (A213)
$ws.Range("B213").Value = 11
$ws.Range("D213").Value = 10
$ws.Range("E213").Value = 0
$ws.Range("H213").Value = 1
